# Refresh the cryptos price table (columns D = Price, E = Volume(1h)%)
# with the latest scrape, preserving each cell's original text-like storage
# (prices such as "209.42" must stay text, not be coerced to numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "27.622.68"; E = "  -1.21%  " }
    @{ Row = 3; D = "1.617.59"; E = "  -1.12%  " }
    @{ Row = 4; D = $null; E = "  -0.88%  " }
    @{ Row = 5; D = "209.42"; E = "  -1.34%  " }
    @{ Row = 6; D = $null; E = "  -1.46%  " }
    @{ Row = 7; D = "0.991"; E = "  -0.85%  " }
    @{ Row = 8; D = "23.13"; E = "  -0.81%  " }
    @{ Row = 9; D = "0.255"; E = "  -1.26%  " }
    @{ Row = 10; D = "0.0605"; E = "  -1.59%  " }
    @{ Row = 11; D = $null; E = "  -0.93%  " }
    @{ Row = 12; D = "1.848.39"; E = "  -1.06%  " }
    @{ Row = 13; D = "1.623.81"; E = "  -0.59%  " }
    @{ Row = 14; D = $null; E = "  -1.94%  " }
    @{ Row = 15; D = "0.556"; E = "  -1.63%  " }
    @{ Row = 16; D = "64.63"; E = "  -1.19%  " }
    @{ Row = 17; D = "27.659.89"; E = "  -1.07%  " }
    @{ Row = 18; D = "226.82"; E = "  -1.79%  " }
    @{ Row = 19; D = "7.63"; E = "  +1.10%  " }
    @{ Row = 20; D = "0.0₃0714"; E = "  -1.37%  " }
    @{ Row = 21; D = "0.992"; E = "  -0.84%  " }
    @{ Row = 22; D = "4.30"; E = "  -1.66%  " }
    @{ Row = 23; D = $null; E = "  -2.95%  " }
    @{ Row = 24; D = "2.04"; E = "  -0.93%  " }
    @{ Row = 25; D = "154.38"; E = "  -0.33%  " }
    @{ Row = 26; D = "6.89"; E = "  -1.13%  " }
    @{ Row = 27; D = $null; E = "  -1.00%  " }
    @{ Row = 28; D = "15.39"; E = "  -1.79%  " }
    @{ Row = 29; D = "0.992"; E = "  -0.87%  " }
    @{ Row = 30; D = $null; E = "  -1.01%  " }
    @{ Row = 31; D = $null; E = "  -1.25%  " }
    @{ Row = 32; D = "3.37"; E = "  -0.96%  " }
    @{ Row = 33; D = "3.07"; E = "  -0.23%  " }
    @{ Row = 34; D = "1.389.61"; E = "  -1.26%  " }
    @{ Row = 35; D = $null; E = "  +1.41%  " }
    @{ Row = 36; D = "1.00"; E = "  -1.13%  " }
    @{ Row = 37; D = "2.32"; E = "  -1.68%  " }
    @{ Row = 38; D = $null; E = "  -0.11%  " }
    @{ Row = 39; D = "0.555"; E = "  -1.68%  " }
    @{ Row = 40; D = "0.842"; E = "  -3.45%  " }
    @{ Row = 41; D = $null; E = "  -1.20%  " }
    @{ Row = 42; D = $null; E = "  -0.85%  " }
    @{ Row = 43; D = $null; E = "  -0.87%  " }
    @{ Row = 44; D = "65.37"; E = "  -2.13%  " }
    @{ Row = 45; D = "5.38"; E = "  -2.96%  " }
    @{ Row = 46; D = "1.759.05"; E = "  -1.10%  " }
    @{ Row = 47; D = "2.11"; E = "  -4.01%  " }
    @{ Row = 48; D = "87.58"; E = "  -0.49%  " }
    @{ Row = 49; D = "0.101"; E = "  +1.16%  " }
    @{ Row = 50; D = $null; E = "  -0.79%  " }
    @{ Row = 51; D = "7.54"; E = "  +0.48%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $priceCell = $ws.Cells.Item($u.Row, 4)
        # Force text storage so numeric-looking prices (e.g. "209.42")
        # are not auto-converted to numbers by the COM Value setter,
        # then restore the default "Normal" style so formatting is untouched.
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.D
        $priceCell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
